$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Enkosi ngokuthatha inxaxheba kolu dliwano-ndlebe. Kuya kuthatha malunga nemizuzu engamashumi amathathu. I will need to record this conversation, with your permission, so that we can listen to your contributions at a later stage. Nangona kunjalo, zonke iimpendulo zakho ziya kuba yimfihlo ngokupheleleyo kwaye ziya kujongwa kuphela liqela lophando. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enkosi ngokuthatha inxaxheba kolu dliwano-ndlebe. Kuya kuthatha malunga nemizuzu engamashumi amathathu. Kuya kufuneka siyirekhode lengxoxo, ngemvume yakho, ukuze simamele igalelo lakho emva kwexesha. Nangona kunjalo, zonke iimpendulo zakho ziya kuba yimfihlo ngokupheleleyo kwaye ziya kujongwa kuphela liqela lophando. ",
    2
)

$d.Content.Find.Execute(
    "Buza: Ziye zathini izigulane malunga neepowusta? Did they mention seeing the posters or say anything else about the posters? ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Buza: Ziye zathini izigulane malunga neepowusta? Ingaba baye bakhankanya ngokubona iipowusta okanye bathetha enye into malunga neepowusta? ",
    2
)

$d.Content.Find.Execute(
    "Could you tell me what you think your overall sense of the patients’ overall interest in ParentText was? Batheni okanye benze ntoni ekwenze ucinge ngolu hlobo?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ingaba ungandixelela ukuba ucinga ukuba yintoni na imvakalelo yakho iyonke yomdla uwonke wezigulana kwi-ParentText? Batheni okanye benze ntoni ekwenze ucinge ngolu hlobo?",
    2
)

$d.Content.Find.Execute(
    "What is your impression of nurses being the encouragers for this intervention?  Ingaba ucinga ukuba a) kuyenzeka kwaye b) kusengqiqweni ukuqhubeka nokucela abongikazi ukuba bakhuthaze inkqubo ngexesha lothethwano lwabo lwesiqhelo nezigulane? Kutheni/kutheni kungenjalo ku-a) kunye no-b)?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Uthini umbono wakho ngokuba abongikazi ibengabakhuthazi kolu ngenelelo?  Ingaba ucinga ukuba a) kuyenzeka kwaye b) kusengqiqweni ukuqhubeka nokucela abongikazi ukuba bakhuthaze inkqubo ngexesha lothethwano lwabo lwesiqhelo nezigulane? Kutheni/kutheni kungenjalo ku-a) kunye no-b)?",
    2
)
